$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage so that
# numeric-looking strings (e.g. "1.11") are not auto-converted to
# numbers by Excel, and the cells original style/number-format
# is preserved (no stray formatting diffs).
function Set-TextValue($cell, [string]$value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range('D2') '55.675.62'
$ws.Range('E2').Value = '  +2.67%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.499.80'
$ws.Range('E3').Value = '  +8.79%  '

# Row 4
$ws.Range('E4').Value = '  +0.11%  '

# Row 5
Set-TextValue $ws.Range('D5') '480.24'
$ws.Range('E5').Value = '  +6.79%  '

# Row 6
Set-TextValue $ws.Range('D6') '138.45'
$ws.Range('E6').Value = '  +6.88%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.999'
$ws.Range('E7').Value = '  +0.34%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.512'
$ws.Range('E8').Value = '  +8.26%  '

# Row 9
Set-TextValue $ws.Range('D9') '2.490.61'
$ws.Range('E9').Value = '  +8.97%  '

# Row 10
$ws.Range('E10').Value = '  +7.38%  '

# Row 11
Set-TextValue $ws.Range('D11') '5.45'
$ws.Range('E11').Value = '  +1.38%  '

# Row 12
$ws.Range('E12').Value = '  +5.15%  '

# Row 13
$ws.Range('E13').Value = '  -0.02%  '

# Row 14
Set-TextValue $ws.Range('D14') '2.928.27'
$ws.Range('E14').Value = '  +9.11%  '

# Row 15
Set-TextValue $ws.Range('D15') '55.701.47'
$ws.Range('E15').Value = '  +2.56%  '

# Row 16
$ws.Range('E16').Value = '  +14.23%  '

# Row 17
$ws.Range('E17').Value = '  +8.70%  '

# Row 18
Set-TextValue $ws.Range('D18') '2.495.11'
$ws.Range('E18').Value = '  +8.74%  '

# Row 19
Set-TextValue $ws.Range('D19') '4.33'
$ws.Range('E19').Value = '  +6.55%  '

# Row 20
Set-TextValue $ws.Range('D20') '320.32'
$ws.Range('E20').Value = '  +6.20%  '

# Row 21
Set-TextValue $ws.Range('D21') '9.95'
$ws.Range('E21').Value = '  +5.48%  '

# Row 22
Set-TextValue $ws.Range('D22') '0.998'
$ws.Range('E22').Value = '  -0.09%  '

# Row 23
$ws.Range('E23').Value = '  +7.33%  '

# Row 24
Set-TextValue $ws.Range('D24') '57.81'
$ws.Range('E24').Value = '  +4.01%  '

# Row 25
$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range('D25') '1.01'
$ws.Range('E25').Value = '  +1.21%  '

# Row 26
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D26') '0.404'
$ws.Range('E26').Value = '  +8.86%  '

# Row 27
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D27') '0.163'
$ws.Range('E27').Value = '  +3.19%  '

# Row 28
Set-TextValue $ws.Range('D28') '2.609.85'
$ws.Range('E28').Value = '  +9.13%  '

# Row 29
Set-TextValue $ws.Range('D29') '7.38'
$ws.Range('E29').Value = '  +7.95%  '

# Row 30
Set-TextValue $ws.Range('D30') '0.0₃0772'
$ws.Range('E30').Value = '  +8.77%  '

# Row 31
$ws.Range('E31').Value = '  +0.47%  '

# Row 32
Set-TextValue $ws.Range('D32') '148.10'
$ws.Range('E32').Value = '  +2.05%  '

# Row 33
$ws.Range('E33').Value = '  +6.63%  '

# Row 34
$ws.Range('E34').Value = '  +9.42%  '

# Row 35
Set-TextValue $ws.Range('D35') '5.17'
$ws.Range('E35').Value = '  +9.60%  '

# Row 36
$ws.Range('E36').Value = '  +1.90%  '

# Row 37
Set-TextValue $ws.Range('D37') '1.11'
$ws.Range('E37').Value = '  +9.35%  '

# Row 38
Set-TextValue $ws.Range('D38') '0.838'
$ws.Range('E38').Value = '  -0.82%  '

# Row 39
Set-TextValue $ws.Range('D39') '34.47'
$ws.Range('E39').Value = '  +4.58%  '

# Row 40
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D40') '0.998'
$ws.Range('E40').Value = '  +0.51%  '

# Row 41
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D41') '0.612'
$ws.Range('E41').Value = '  +19.69%  '

# Row 42
$ws.Range('E42').Value = '  +9.92%  '

# Row 43
Set-TextValue $ws.Range('D43') '1.32'
$ws.Range('E43').Value = '  +6.08%  '

# Row 44
$ws.Range('E44').Value = '  +6.70%  '

# Row 45
Set-TextValue $ws.Range('D45') '10.17'
$ws.Range('E45').Value = '  -1.20%  '

# Row 46
Set-TextValue $ws.Range('D46') '1.963.68'
$ws.Range('E46').Value = '  +1.51%  '

# Row 47
Set-TextValue $ws.Range('D47') '0.0902'
$ws.Range('E47').Value = '  +10.55%  '

# Row 48
$ws.Range('E48').Value = '  +7.88%  '

# Row 49
Set-TextValue $ws.Range('D49') '250.28'
$ws.Range('E49').Value = '  +31.60%  '

# Row 50
Set-TextValue $ws.Range('D50') '17.47'
$ws.Range('E50').Value = '  +8.00%  '

# Row 51
Set-TextValue $ws.Range('D51') '4.39'
$ws.Range('E51').Value = '  +8.82%  '
